# The "DailyProductionReport" sheet is a daily report template that, in this
# revision, gets reset back to a blank template: the five data rows (10-14,
# columns A:O) are emptied out (values + formulas removed, formatting kept)
# and the active selection is moved to O10. The totals row (15) keeps its
# SUM formulas, which now simply evaluate against the now-empty range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the five detail rows (keeps cell styles/formatting intact,
# removes literal values and the per-row J10:J14 formulas).
$ws.Range("A10:O14").ClearContents()

# Move/restore the visible selection to O10, matching the saved view state.
$ws.Range("O10").Select()
